$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- multiple inheritance section ---
$ws.Range("A40").Value = "multiple inheritance:"

$ws.Range("B41").Value = "create a FooBar"
$ws.Range("C41").Formula = "=_xll.clFooBar(,""foobar"")"

$ws.Range("B42").Value = "use it as a Foo"
$ws.Range("C42").Formula = "=_xll.clFunctionUsingFoo(,C41)"

$ws.Range("B43").Value = "use it as a Bar"
$ws.Range("C43").Formula = "=_xll.clFunctionUsingBar(,C41)"

# --- diamond inheritance section ---
$ws.Range("A44").Value = "diamond inheritance:"

$ws.Range("B45").Value = "create a Foo2"
$ws.Range("C45").Formula = "=_xll.clFoo2(,""foo2"")"

$ws.Range("B46").Value = "create a Bar2"
$ws.Range("C46").Formula = "=_xll.clBar2(,""bar2"")"

$ws.Range("B47").Value = "create a FooBar2"
$ws.Range("C47").Formula = "=_xll.clFooBar2(,""foobar2"")"

# NOTE: the shared-string table indexes strings in first-use order, and the
# target workbook's table has the "use fooX as a FooX/Bar2" strings (rows
# 52-55) ahead of the "call FooX::f() ..." strings (rows 48-51). Assign the
# cell *values* in that same order so the shared-string indices line up,
# then fill in the formulas for rows 48-51 afterwards.

$ws.Range("B52").Value = "use foo2 as a Foo2"
$ws.Range("B53").Value = "use foobar2 as a Foo2"
$ws.Range("B54").Value = "use bar2 as a Bar2"
$ws.Range("B55").Value = "use foobar2 as a Bar2"

$ws.Range("B48").Value = "call Foo2::f() on foo2"
$ws.Range("C48").Formula = "=_xll.clFoo2F(,C45)"

$ws.Range("B49").Value = "call Foo2::f() on foobar2"
$ws.Range("C49").Formula = "=_xll.clFoo2F(,C47)"

$ws.Range("B50").Value = "call Bar2::f() on bar2"
$ws.Range("C50").Formula = "=_xll.clBar2F(,C46)"

$ws.Range("B51").Value = "call Bar2::f() on foobar2"
$ws.Range("C51").Formula = "=_xll.clBar2F(,C47)"

$ws.Range("C52").Formula = "=_xll.clFunctionUsingFoo2(,C45)"
$ws.Range("C53").Formula = "=_xll.clFunctionUsingFoo2(,C47)"
$ws.Range("C54").Formula = "=_xll.clFunctionUsingBar2(,C46)"
$ws.Range("C55").Formula = "=_xll.clFunctionUsingBar2(,C47)"
